$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1947.8572
$ws.Range("I2").Value = 222.85715
$ws.Range("K2").Value = 222.85715
$ws.Range("M2").Value = -109.85715
$ws.Range("H6").Value = 437.8889
$ws.Range("I6").Value = 277.57144
$ws.Range("J6").Value = 999
$ws.Range("K6").Value = 832.71432
$ws.Range("L6").Value = 2997
$ws.Range("M6").Value = -720.71432
$ws.Range("N6").Value = -3221
$ws.Range("H106").Value = 1100
$ws.Range("I106").Value = 1100
$ws.Range("K106").Value = 1100
$ws.Range("M106").Value = -469
$ws.Range("H116").Value = 7980.273
$ws.Range("J116").Value = 5964.1665
$ws.Range("L116").Value = 5964.1665
$ws.Range("N116").Value = -12848.1665
$ws.Range("H137").Value = 1433.7273
$ws.Range("I137").Value = 1433.7273
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 4301.1819
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -1751.1819
$ws.Range("N137").ClearContents()
$ws.Range("H138").Value = 3405.3125
$ws.Range("I138").Value = 2570.8572
$ws.Range("K138").Value = 7712.571599999999
$ws.Range("M138").Value = -2572.571599999999
$ws.Range("H141").Value = 3355
$ws.Range("I141").Value = 3297.2
$ws.Range("K141").Value = 9891.599999999999
$ws.Range("M141").Value = -4711.599999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5139806.5
$ws.Range("I32").Value = 5006935.5
$ws.Range("K32").Value = 5006935.5
$ws.Range("M32").Value = -5006648.5
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("H97").Value = 923.63635
$ws.Range("I97").Value = 875
$ws.Range("K97").Value = 875
$ws.Range("M97").Value = -379
$ws.Range("H132").Value = 2389.4
$ws.Range("I132").Value = 2389.4
$ws.Range("K132").Value = 7168.200000000001
$ws.Range("M132").Value = -4638.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 125
$ws.Range("I80").Value = 111
$ws.Range("J80").Value = 167
$ws.Range("K80").Value = 111
$ws.Range("L80").Value = 167
$ws.Range("M80").Value = 887
$ws.Range("N80").Value = -2163
$ws.Range("H83").Value = 125
$ws.Range("I83").Value = 111
$ws.Range("J83").Value = 167
$ws.Range("K83").Value = 555
$ws.Range("L83").Value = 835
$ws.Range("M83").Value = 4437
$ws.Range("N83").Value = -10819
$ws.Range("H134").Value = 1688.3334
$ws.Range("I134").Value = 1688.3334
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5065.0002
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2530.0002
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2354.3845
$ws.Range("I31").Value = 2180.25
$ws.Range("K31").Value = 2180.25
$ws.Range("M31").Value = -1885.25
$ws.Range("H34").Value = 2354.3845
$ws.Range("I34").Value = 2180.25
$ws.Range("K34").Value = 2180.25
$ws.Range("M34").Value = -1978.25
$ws.Range("H86").Value = 7227.75
$ws.Range("I86").Value = 7340.636
$ws.Range("K86").Value = 7340.636
$ws.Range("M86").Value = -6217.636
$ws.Range("H89").Value = 7227.75
$ws.Range("I89").Value = 7340.636
$ws.Range("K89").Value = 36703.18
$ws.Range("M89").Value = -31087.18
$ws.Range("H105").Value = 3285.7144
$ws.Range("I105").Value = 2375
$ws.Range("K105").Value = 2375
$ws.Range("M105").Value = -628
$ws.Range("H134").Value = 3284.5625
$ws.Range("I134").Value = 2712.75
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 8138.25
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -5603.25
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 500
$ws.Range("I43").Value = 600
$ws.Range("K43").Value = 1800
$ws.Range("M43").Value = -1686
$ws.Range("H97").Value = 2744
$ws.Range("I97").Value = 2744
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 8232
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -7736
$ws.Range("N97").ClearContents()
$ws.Range("H106").Value = 18642.715
$ws.Range("I106").Value = 16833.334
$ws.Range("K106").Value = 50500.00199999999
$ws.Range("M106").Value = -49554.00199999999
$ws.Range("H109").Value = 1513.5
$ws.Range("I109").Value = 1513.5
$ws.Range("K109").Value = 4540.5
$ws.Range("M109").Value = -3500.5
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 13005.4
$ws.Range("I24").Value = 14006
$ws.Range("J24").Value = 12755.25
$ws.Range("K24").Value = 14006
$ws.Range("L24").Value = 12755.25
$ws.Range("M24").Value = -13833
$ws.Range("N24").Value = -13101.25
$ws.Range("H99").Value = 7798.3335
$ws.Range("I99").Value = 4454.6
$ws.Range("J99").Value = 11978
$ws.Range("K99").Value = 4454.6
$ws.Range("L99").Value = 11978
$ws.Range("M99").Value = -2208.6
$ws.Range("N99").Value = -16470
$ws.Range("H134").Value = 59000
$ws.Range("J134").Value = 59000
$ws.Range("L134").Value = 177000
$ws.Range("N134").Value = -182070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 3937.3333
$ws.Range("J20").Value = 3468.5
$ws.Range("L20").Value = 3468.5
$ws.Range("N20").Value = -3920.5
$ws.Range("H23").Value = 10000
$ws.Range("I23").Value = 10000
$ws.Range("K23").Value = 10000
$ws.Range("M23").Value = -9770

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H63").Value = 40499.332
$ws.Range("J63").Value = 40499.332
$ws.Range("L63").Value = 40499.332
$ws.Range("N63").Value = -41747.332
$ws.Range("H66").Value = 40499.332
$ws.Range("J66").Value = 40499.332
$ws.Range("L66").Value = 121497.996
$ws.Range("N66").Value = -127737.996
$ws.Range("H69").Value = 7300
$ws.Range("I69").Value = 7300
$ws.Range("K69").Value = 7300
$ws.Range("M69").Value = -6551
$ws.Range("H72").Value = 7300
$ws.Range("I72").Value = 7300
$ws.Range("K72").Value = 21900
$ws.Range("M72").Value = -18156
$ws.Range("H82").Value = 40000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 40000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 40000
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -40766
$ws.Range("H85").Value = 40000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 40000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 40000
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -42652
$ws.Range("H112").Value = 63462.332
$ws.Range("J112").Value = 63462.332
$ws.Range("L112").Value = 63462.332
$ws.Range("N112").Value = -66416.33199999999
$ws.Range("H113").Value = 339
$ws.Range("I113").Value = 348.75
$ws.Range("J113").Value = 329.25
$ws.Range("K113").Value = 1046.25
$ws.Range("L113").Value = 987.75
$ws.Range("M113").Value = 1123.75
$ws.Range("N113").Value = -5327.75
$ws.Range("H132").Value = 6890.5
$ws.Range("I132").Value = 6890.5
$ws.Range("K132").Value = 20671.5
$ws.Range("M132").Value = -18141.5
$ws.Range("H136").Value = 2121.1428
$ws.Range("I136").Value = 1641.3334
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 4924.0002
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -2374.0002
$ws.Range("N136").Value = -20100
